$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A.
# This shifts the existing Code/Description/Definition columns from A/B/C to B/C/D.
$ws.Columns.Item(1).Insert()

# New header cells: A1 = "Version", E1 = "Guide.for.Use"
# (B1/C1/D1 keep their original "Code"/"Description"/"Definition" values.)
$ws.Range("A1").Value = "Version"
$ws.Range("E1").Value = "Guide.for.Use"

# Fill the new "Version" column (A2:A7) with the text value "1.0".
# Applying a text number format first keeps Excel from auto-converting the
# numeric-looking string "1.0" into the number 1; ClearFormats afterwards
# removes the temporary formatting again so no extra style is left applied
# to the cells (the text type set on the cell itself is preserved).
$verRange = $ws.Range("A2:A7")
$verRange.NumberFormat = "@"
$ws.Range("A2").Value = "1.0"
$ws.Range("A3").Value = "1.0"
$ws.Range("A4").Value = "1.0"
$ws.Range("A5").Value = "1.0"
$ws.Range("A6").Value = "1.0"
$ws.Range("A7").Value = "1.0"
$verRange.ClearFormats()

# Create the new, empty "Guide.for.Use" data cells E2:E7.
# Briefly formatting + clearing makes the engine materialize an (empty)
# cell entry at each address, matching the target's bare <c r="E2"/> cells.
$guideRange = $ws.Range("E2:E7")
$guideRange.NumberFormat = "@"
$guideRange.ClearFormats()
